$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer text (A59)
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-08 for illustrative purposes only and are subject to change."
$ws.Rows(59).EntireRow.AutoFit()

# Update Weight (D) and Percent Change (E) values for rows 2-56
$ws.Range("D2").Value = 0.01863656539172915
$ws.Range("E2").Value = -0.007854984894259731
$ws.Range("D3").Value = 0.01816736686726467
$ws.Range("E3").Value = -0.01570247933884306
$ws.Range("D4").Value = 0.01987750164923281
$ws.Range("E4").Value = 0.005268524813051956
$ws.Range("D5").Value = 0.02003965665928773
$ws.Range("E5").Value = -0.02000449539222293
$ws.Range("D6").Value = 0.02014438176994821
$ws.Range("E6").Value = -0.02163340600368946
$ws.Range("D7").Value = 0.008012034003755464
$ws.Range("E7").Value = -0.01171234481143124
$ws.Range("D8").Value = 0.01962338372838285
$ws.Range("E8").Value = -0.01451825780906302
$ws.Range("D9").Value = 0.02288750402337734
$ws.Range("E9").Value = 0.002952029520295163
$ws.Range("D10").Value = 0.02314424945596431
$ws.Range("E10").Value = -0.02067824648469807
$ws.Range("D11").Value = 0.01958021746413212
$ws.Range("E11").Value = -0.03180354267310781
$ws.Range("D12").Value = 0.01417917940931659
$ws.Range("E12").Value = -0.0003970880211781047
$ws.Range("D13").Value = 0.01510068531136483
$ws.Range("E13").Value = -0.02326621923937378
$ws.Range("D14").Value = 0.00897783224651315
$ws.Range("E14").Value = -0.00886361735931096
$ws.Range("D15").Value = 0.01466601979830094
$ws.Range("E15").Value = -0.03194103194103204
$ws.Range("D16").Value = 0.02328932563972873
$ws.Range("E16").Value = 0.005673255917028497
$ws.Range("D17").Value = 0.02486677107897831
$ws.Range("E17").Value = -0.01648351648351665
$ws.Range("D18").Value = 0.02318553892611718
$ws.Range("E18").Value = -0.02214703168255916
$ws.Range("D19").Value = 0.01800427346016082
$ws.Range("E19").Value = -0.03577571379428968
$ws.Range("D20").Value = 0.01970164604226353
$ws.Range("E20").Value = -0.01755656108597292
$ws.Range("D21").Value = 0.0282842254517678
$ws.Range("E21").Value = -0.01635645798082352
$ws.Range("D22").Value = 0.01815291555271117
$ws.Range("E22").Value = -0.004001116590676479
$ws.Range("D23").Value = 0.0204998465720825
$ws.Range("E23").Value = -0.02819789797487815
$ws.Range("D24").Value = 0.01870131478810525
$ws.Range("E24").Value = -0.02051282051282055
$ws.Range("D25").Value = 0.01936870276930353
$ws.Range("E25").Value = 0.01181190104747043
$ws.Range("D26").Value = 0.01815798289677538
$ws.Range("E26").Value = -0.01534883720930236
$ws.Range("D27").Value = 0.02279291360084531
$ws.Range("E27").Value = -0.02718903874973233
$ws.Range("D28").Value = 0.02227942273567138
$ws.Range("E28").Value = -0.01415213545615357
$ws.Range("D29").Value = 0.01976939830919619
$ws.Range("E29").Value = -0.02412280701754377
$ws.Range("D30").Value = 0.02162235712201132
$ws.Range("E30").Value = -0.01239486498450615
$ws.Range("D31").Value = 0.01995125965727863
$ws.Range("E31").Value = -0.01179624664879353
$ws.Range("D32").Value = 0.02096641758480998
$ws.Range("E32").Value = -0.02155504234026162
$ws.Range("D33").Value = 0.01766438604903875
$ws.Range("E33").Value = -0.01823204419889513
$ws.Range("D34").Value = 0.01860860115967107
$ws.Range("E34").Value = -0.01140684410646398
$ws.Range("D35").Value = 0.02125356708178224
$ws.Range("E35").Value = -0.0008477270318957286
$ws.Range("D36").Value = 0.01773570422475735
$ws.Range("E36").Value = -0.01190476190476186
$ws.Range("D37").Value = 0.02041351404358103
$ws.Range("E37").Value = -0.01072006472491893
$ws.Range("D38").Value = 0.01890626070359134
$ws.Range("E38").Value = 0.006432591798445353
$ws.Range("D39").Value = 0.01864219577402273
$ws.Range("E39").Value = -0.005919661733615134
$ws.Range("D40").Value = 0.016884390421969
$ws.Range("E40").Value = -0.01980792316926772
$ws.Range("D41").Value = 0.01327268786005122
$ws.Range("E41").Value = -0.008936651583710487
$ws.Range("D42").Value = 0.01387326197136575
$ws.Range("E42").Value = 0.0345238095238094
$ws.Range("D43").Value = 0.01662595587469396
$ws.Range("E43").Value = -0.01535213970447136
$ws.Range("D44").Value = 0.01277984172995373
$ws.Range("E44").Value = -0.009736540664375903
$ws.Range("D45").Value = 0.01519677716917516
$ws.Range("E45").Value = -0.01952526799387455
$ws.Range("D46").Value = 0.02102760107240015
$ws.Range("E46").Value = -0.04837558014994647
$ws.Range("D47").Value = 0.01359211821550663
$ws.Range("E47").Value = -0.01691474966170503
$ws.Range("D48").Value = 0.0202224564044191
$ws.Range("E48").Value = -0.03136890951276095
$ws.Range("D49").Value = 0.01814784820864695
$ws.Range("E49").Value = -0.0135372714486639
$ws.Range("D50").Value = 0.01771749932200813
$ws.Range("E50").Value = -0.002129169623846683
$ws.Range("D51").Value = 0.0189755144058023
$ws.Range("E51").Value = -0.00382766601388651
$ws.Range("D52").Value = 0.006486200402196975
$ws.Range("E52").Value = -0.01458333333333339
$ws.Range("D53").Value = 0.02123198394965687
$ws.Range("E53").Value = -0.00440205429200291
$ws.Range("D54").Value = 0.01756360220598378
$ws.Range("E54").Value = -0.01084598698481565
$ws.Range("D55").Value = 0.02054714178334852
$ws.Range("E55").Value = -0.007773109243697451
$ws.Range("D56").Value = 0.9999999999999998
$ws.Range("E56").Value = -0.01357372563334791

$ws.Protect()
